# ---------------------------------------------------------------------------
# Adds a new "Payment Page" worksheet (copied/adapted from "Checkout"),
# updates the Cover Page roll-up counts/formula, tweaks a couple of view
# selections, adds two new rows of test cases to "Checkout", removes the
# stale tabSelected flag from "Product Deatil Page", and finally leaves
# "Payment Page" as the active/selected sheet - matching the target diff.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Checkout sheet: two additional test-case rows (30 & 31)
# ---------------------------------------------------------------------------
$checkout = $wb.Worksheets.Item("Checkout")

$checkout.Range("A30").Value = "TC_14"
$checkout.Range("B30").Value = "1. User is on Cart page with more than 2 products and click on secure checkout.`n2. User is on Checkout page.`n3. Enter information.`n4. Click on ""Proceed to payment"" button.`n"
$checkout.Range("C30").Value = 'Verify application displays payment page when user click on "Proceed to payment" button.'
$checkout.Rows.Item(30).RowHeight = 120

$checkout.Range("A31").Value = "TC_15"
$checkout.Range("B31").Value = "1. User is on Cart page with more than 2 products and click on secure checkout.`n2. User is on Checkout page.`n3. Click on website logo.`n"
$checkout.Range("C31").Value = 'Verify application display home page when user click on website logo display on topleft corner.'
$checkout.Rows.Item(31).RowHeight = 90

$checkout.Range("B30:C31").WrapText = $true
$checkout.Range("B30:C31").VerticalAlignment = -4160

# formatting tweak: enable wrap text on C32:C37 (was plain before)
$checkout.Range("C32:C37").WrapText = $true
$checkout.Range("C32:C37").VerticalAlignment = -4160

# updated scroll/selection on the Checkout sheet
$checkout.Range("C20").Select()

# ---------------------------------------------------------------------------
# 2. New "Payment Page" sheet - created as a copy of "Checkout" so it keeps
#    the same summary widgets / layout / drawing, then adapted.
# ---------------------------------------------------------------------------
$checkout.Copy([System.Reflection.Missing]::Value, $checkout)
$payment = $wb.Worksheets.Item($checkout.Index + 1)
$payment.Name = "Payment Page"

# Drop everything below row 26 - Payment Page only needs rows 1-26.
$payment.Range("A27:G55").Clear()

# Row 17
$payment.Range("B17").Value = "1. Launch Ninja specific URL.`n2. User is on Home Page`n3. Click on product that have options.`n4. User is on product detail page.`n5. Enter valid information and click on Check availability button`n6. User is on availability page and click on add to cart.`n7. User is on Cart page and click on secure checkout`n8. User is on Checkout page.`n9. Enter information.`n10. Click on ""Proceed to payment"" button."
$payment.Range("C17").Value = 'Verify application displays payment page when user click on "Proceed to payment" button.'

# Row 18 (header verification reuses existing wording - only Step changes)
$payment.Range("B18").Value = "1. User is on payment page."

# Row 19
$payment.Range("B19").Value = "1. User is on payment page."
$payment.Range("C19").Value = "Verify Application display below mentioned tabs:`n1. Your details `n2. Payment (Highlighted)`n3. Confirm"

# Row 20 (footer verification reuses existing wording - only Step changes)
$payment.Range("B20").Value = "1. User is on payment page."

# Row 21
$payment.Range("B21").ClearContents()
$payment.Range("C21").Value = "Verify application display below mentioned areas on payment page:`n1. Contact Information `n2. "

# Rows 22-26: clear the leftover Checkout-specific wording, keep only the
# step-number column (A) filled in.
$payment.Range("B22:E26").ClearContents()

# normalise wrap/alignment across the body of the table to match the
# intended look (wrap text, top aligned) while keeping the outer D/E/F/G
# columns on row 17 unwrapped (matches the template it was copied from).
$payment.Range("B17:B26").WrapText = $true
$payment.Range("B17:B26").VerticalAlignment = -4160
$payment.Range("C17:C18").WrapText = $true
$payment.Range("C17:C18").VerticalAlignment = -4160
$payment.Range("C19:C20").WrapText = $true
$payment.Range("C21:C26").WrapText = $true
$payment.Range("C21:C26").VerticalAlignment = -4160
$payment.Range("D18:E26").WrapText = $true
$payment.Range("D18:E26").VerticalAlignment = -4160

$payment.Rows.Item(17).RowHeight = 225
$payment.Rows.Item(18).RowHeight = 180
$payment.Rows.Item(19).RowHeight = 75
$payment.Rows.Item(20).RowHeight = 120
$payment.Rows.Item(21).RowHeight = 60

$payment.Range("C21").Select()

# ---------------------------------------------------------------------------
# 3. Cover Page: refresh the roll-up table with a "Payment Page" row and
#    updated totals.
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover Page")
$cover.Range("B29").Formula = "=COUNTIF('Payment Page'!A17:A70,""*"")"
$cover.Range("B29:C29").Select()

# ---------------------------------------------------------------------------
# 4. "Product Deatil Page" is no longer the active tab - it keeps its own
#    scroll position instead.
# ---------------------------------------------------------------------------
$productDetail = $wb.Worksheets.Item("Product Deatil Page")
$productDetail.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 5. Finally, "Payment Page" becomes the active/selected sheet.
# ---------------------------------------------------------------------------
$payment.Activate()
$wb.Windows.Item(1).TabRatio = 725
